$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.043.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.55%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.827.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.45%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.22%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6364"
$ws.Range("D6").Style = "Normal"

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.73"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.40%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2935"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.41%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07330"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.57%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.78"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.49%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07658"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.74%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.827.05"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.23%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.980"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.02%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6631"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.50%  "

$ws.Range("E16").Value = "  -1.89%  "

$ws.Range("E17").Value = "  -0.65%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008654"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.66%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "28.899.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.92%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.075.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.27%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.42%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.95%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("E24").Value = "  -0.16%  "

$ws.Range("E25").Value = "  +0.09%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.73%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.464"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1368"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "17.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.37%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.504"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.090"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.54%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.024"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.20%  "

$ws.Range("E33").Value = "  +1.52%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05291"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.21%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.835"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.75%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7373"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.17%  "

$ws.Range("E37").Value = "  +1.83%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.655"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.88%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.291.27"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.82%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.744"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.96%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01780"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.95%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.288"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.72%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8949"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.93%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.09%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.59"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.28%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.975.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.11%  "

$ws.Range("E47").Value = "  -0.47%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "64.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.11%  "

$ws.Range("E49").Value = "  -5.51%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.727"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.23%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07302"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -15.27%  "
